$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.033.40'
$ws.Range("E2").Value = '  +3.13%  '
$ws.Range("D3").Value = '3.239.48'
$ws.Range("E3").Value = '  +7.18%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.20'
$ws.Range("E5").Value = '  +5.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.71'
$ws.Range("E6").Value = '  +9.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.229.93'
$ws.Range("E9").Value = '  +6.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.11'
$ws.Range("E10").Value = '  +12.14%  '
$ws.Range("E11").Value = '  +7.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.487'
$ws.Range("E12").Value = '  +6.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.65'
$ws.Range("E13").Value = '  +4.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000234'
$ws.Range("E14").Value = '  +7.60%  '
$ws.Range("D15").Value = '3.756.13'
$ws.Range("E15").Value = '  +7.21%  '
$ws.Range("D16").Value = '66.101.09'
$ws.Range("E16").Value = '  +3.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '545.15'
$ws.Range("E17").Value = '  +13.96%  '
$ws.Range("D18").Value = '3.240.01'
$ws.Range("E18").Value = '  +7.31%  '
$ws.Range("E19").Value = '  +3.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.10'
$ws.Range("E20").Value = '  +7.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.51'
$ws.Range("E21").Value = '  +7.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.743'
$ws.Range("E22").Value = '  +9.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.88'
$ws.Range("E23").Value = '  +12.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.46'
$ws.Range("E24").Value = '  +8.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.02'
$ws.Range("E25").Value = '  +4.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.28'
$ws.Range("E27").Value = '  +19.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.98'
$ws.Range("E28").Value = '  +10.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.25'
$ws.Range("E29").Value = '  +7.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.65'
$ws.Range("E30").Value = '  +7.86%  '
$ws.Range("E31").Value = '  +6.74%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.18'
$ws.Range("E33").Value = '  +6.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '567.48'
$ws.Range("E34").Value = '  +10.23%  '
$ws.Range("E35").Value = '  +5.08%  '
$ws.Range("E36").Value = '  +7.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '55.11'
$ws.Range("E37").Value = '  +5.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0452'
$ws.Range("E38").Value = '  +14.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0860'
$ws.Range("E39").Value = '  +8.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.129'
$ws.Range("E40").Value = '  +6.62%  '
$ws.Range("D41").Value = '3.192.83'
$ws.Range("E41").Value = '  +11.18%  '
$ws.Range("E42").Value = '  +9.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.57'
$ws.Range("E43").Value = '  +4.62%  '
$ws.Range("E44").Value = '  +17.83%  '
$ws.Range("E45").Value = '  +12.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.47'
$ws.Range("E46").Value = '  +7.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.999'
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("E48").Value = '  +7.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '126.11'
$ws.Range("E49").Value = '  +5.79%  '
$ws.Range("E50").Value = '  +4.49%  '
$ws.Range("E51").Value = '  +9.83%  '
